# "Generate Report for Handback"
#
# The f8e9b7b9-... file has now been handed back (translations are in sync
# with en-US) while 67db45b4-... is still only "Ready for handoff". This
# script updates the localization-status report:
#   - Overview sheet: the two file rows swap order (handed-back file first)
#     and pick up the new status/date.
#   - zh-cn / de-de sheets: same row re-ordering, the handed-back row gets
#     its Status, Latest Target File, Latest Handback File and Latest
#     Handback DateTime populated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$urlOverview67dbMd = "https://github.com/OpenLocalizationTest/oltest/blob/53766a93dde60e068b355cfd8de756851ddad303/e2e/67db45b4-2ff1-4d56-a21d-cd44083119e0.md"
$urlOverviewF8e9Md = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTest/oltest/ci/f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"

# Row 2 now holds the handed-back file (f8e9...), row 3 the still-pending one (67db...)
$ws1.Range("A2").Value = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-30-20 12:30:45"

$ws1.Range("A3").Value = "67db45b4-2ff1-4d56-a21d-cd44083119e0.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-30-20 12:30:23"

foreach ($hl in $ws1.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.Address = $urlOverviewF8e9Md
        $hl.TextToDisplay = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
    } elseif ($addr -eq '$A$3') {
        $hl.Address = $urlOverview67dbMd
        $hl.TextToDisplay = "67db45b4-2ff1-4d56-a21d-cd44083119e0.md"
    }
}

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$url67dbMd        = "https://github.com/OpenLocalizationTest/oltest/blob/53766a93dde60e068b355cfd8de756851ddad303/e2e/67db45b4-2ff1-4d56-a21d-cd44083119e0.md"
$url67dbZhXlf     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fed196a0930cf8c574edb31ceb954ec9f3391c3d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/67db45b4-2ff1-4d56-a21d-cd44083119e0.94bbdd89e87978290c262512ecf23b43ef0b6516.zh-cn.xlf"
$urlF8e9MdBack    = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTest/oltest/ci/f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
$urlF8e9ZhXlfBack = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTest/oltest/ci/ht/f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.zh-cn.xlf"

# Row 2: the handed-back file (f8e9...)
$ws2.Range("A2").Value = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-20 12:30:42"
$ws2.Range("H2").Value = "2016-03-20 12:31:02"
$ws2.Range("I2").Value = "Include"

# Row 3: still only handed off (67db...)
$ws2.Range("A3").Value = "67db45b4-2ff1-4d56-a21d-cd44083119e0.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "67db45b4-2ff1-4d56-a21d-cd44083119e0.94bbdd89e87978290c262512ecf23b43ef0b6516.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-20 12:30:20"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"

foreach ($hl in $ws2.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.Address = $urlF8e9MdBack
        $hl.TextToDisplay = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
    } elseif ($addr -eq '$B$2') {
        $hl.Address = $urlF8e9MdBack
        $hl.TextToDisplay = ".md"
    } elseif ($addr -eq '$D$2') {
        $hl.Address = $urlF8e9ZhXlfBack
        $hl.TextToDisplay = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.zh-cn.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.Address = $url67dbMd
        $hl.TextToDisplay = "67db45b4-2ff1-4d56-a21d-cd44083119e0.md"
    } elseif ($addr -eq '$B$3') {
        $hl.Address = $url67dbMd
        $hl.TextToDisplay = ".md"
    } elseif ($addr -eq '$D$3') {
        $hl.Address = $url67dbZhXlf
        $hl.TextToDisplay = "67db45b4-2ff1-4d56-a21d-cd44083119e0.94bbdd89e87978290c262512ecf23b43ef0b6516.zh-cn.xlf"
    }
}

# New columns for row 2: the handback produced a target file + handback file
$ws2.Range("F2").Value = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
$ws2.Hyperlinks.Add($ws2.Range("F2"), $urlF8e9MdBack, "", "", "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md") | Out-Null

$ws2.Range("G2").Value = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.zh-cn.xlf"
$ws2.Hyperlinks.Add($ws2.Range("G2"), $urlF8e9ZhXlfBack, "", "", "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$url67dbMdDe       = "https://github.com/OpenLocalizationTest/oltest/blob/53766a93dde60e068b355cfd8de756851ddad303/e2e/67db45b4-2ff1-4d56-a21d-cd44083119e0.md"
$url67dbDeXlf      = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/250baa0b2da41c9bb4d91f47c56963d02d639518/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/67db45b4-2ff1-4d56-a21d-cd44083119e0.94bbdd89e87978290c262512ecf23b43ef0b6516.de-de.xlf"
$urlF8e9MdBackDe   = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTest/oltest/ci/f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
$urlF8e9DeXlfBack  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTest/oltest/ci/ht/f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.de-de.xlf"

# Row 2: the handed-back file (f8e9...)
$ws3.Range("A2").Value = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-20 12:30:45"
$ws3.Range("H2").Value = "2016-03-20 12:31:09"
$ws3.Range("I2").Value = "Include"

# Row 3: still only handed off (67db...)
$ws3.Range("A3").Value = "67db45b4-2ff1-4d56-a21d-cd44083119e0.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "67db45b4-2ff1-4d56-a21d-cd44083119e0.94bbdd89e87978290c262512ecf23b43ef0b6516.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-20 12:30:23"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"

foreach ($hl in $ws3.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.Address = $urlF8e9MdBackDe
        $hl.TextToDisplay = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
    } elseif ($addr -eq '$B$2') {
        $hl.Address = $urlF8e9MdBackDe
        $hl.TextToDisplay = ".md"
    } elseif ($addr -eq '$D$2') {
        $hl.Address = $urlF8e9DeXlfBack
        $hl.TextToDisplay = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.de-de.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.Address = $url67dbMdDe
        $hl.TextToDisplay = "67db45b4-2ff1-4d56-a21d-cd44083119e0.md"
    } elseif ($addr -eq '$B$3') {
        $hl.Address = $url67dbMdDe
        $hl.TextToDisplay = ".md"
    } elseif ($addr -eq '$D$3') {
        $hl.Address = $url67dbDeXlf
        $hl.TextToDisplay = "67db45b4-2ff1-4d56-a21d-cd44083119e0.94bbdd89e87978290c262512ecf23b43ef0b6516.de-de.xlf"
    }
}

# New columns for row 2: the handback produced a target file + handback file
$ws3.Range("F2").Value = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md"
$ws3.Hyperlinks.Add($ws3.Range("F2"), $urlF8e9MdBackDe, "", "", "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.md") | Out-Null

$ws3.Range("G2").Value = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.de-de.xlf"
$ws3.Hyperlinks.Add($ws3.Range("G2"), $urlF8e9DeXlfBack, "", "", "f8e9b7b9-105a-4e28-ae9c-568a30c60d10.e7cbd52349ff8910dbdaba77df21076dc0b9e293.de-de.xlf") | Out-Null
